$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.501.30'

$ws.Range('D3').Value = '1.914.15'
$ws.Range('E3').Value = '  -0.13%  '

$ws.Range('D4').Value = "'1.000"
$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').Value = "'244.66"
$ws.Range('E5').Value = '  +1.48%  '

$ws.Range('D6').Value = "'1.000"
$ws.Range('E6').Value = '  +0.01%  '

$ws.Range('D7').Value = "'0.4839"
$ws.Range('E7').Value = '  +3.14%  '

$ws.Range('D8').Value = "'0.2897"
$ws.Range('E8').Value = '  +1.69%  '

$ws.Range('D9').Value = "'0.06724"
$ws.Range('E9').Value = '  -1.24%  '

$ws.Range('D10').Value = "'109.78"
$ws.Range('E10').Value = '  +2.35%  '

$ws.Range('E11').Value = '  +4.61%  '

$ws.Range('D12').Value = '1.918.75'
$ws.Range('E12').Value = '  +0.19%  '

$ws.Range('D13').Value = "'0.07553"

$ws.Range('D14').Value = "'5.281"
$ws.Range('E14').Value = '  +1.81%  '

$ws.Range('D15').Value = "'0.6730"
$ws.Range('E15').Value = '  +2.88%  '

$ws.Range('D16').Value = "'282.15"
$ws.Range('E16').Value = '  -2.22%  '

$ws.Range('D17').Value = '30.508.01'
$ws.Range('E17').Value = '  +0.32%  '

$ws.Range('D18').Value = "'0.9998"
$ws.Range('E18').Value = '  -0.06%  '

$ws.Range('D19').Value = "'0.000007572"
$ws.Range('E19').Value = '  -0.74%  '

$ws.Range('D20').Value = "'12.88"
$ws.Range('E20').Value = '  -0.54%  '

$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.167.37'
$ws.Range('E21').Value = '  +0.55%  '

$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = "'5.506"
$ws.Range('E22').Value = '  +5.51%  '

$ws.Range('D23').Value = "'1.0000"
$ws.Range('E23').Value = '  +0.07%  '

$ws.Range('D24').Value = "'6.459"
$ws.Range('E24').Value = '  +4.21%  '

$ws.Range('D25').Value = "'9.465"
$ws.Range('E25').Value = '  +2.20%  '

$ws.Range('D26').Value = "'164.28"
$ws.Range('E26').Value = '  -2.10%  '

$ws.Range('D27').Value = "'20.27"
$ws.Range('E27').Value = '  -6.66%  '

$ws.Range('D28').Value = "'2.123"
$ws.Range('E28').Value = '  +4.17%  '

$ws.Range('E29').Value = '  -1.26%  '

$ws.Range('D30').Value = "'1.401"

$ws.Range('D31').Value = "'4.151"
$ws.Range('E31').Value = '  +0.23%  '

$ws.Range('D32').Value = "'4.044"
$ws.Range('E32').Value = '  +2.75%  '

$ws.Range('D33').Value = "'0.04995"
$ws.Range('E33').Value = '  -0.58%  '

$ws.Range('D34').Value = "'0.7310"
$ws.Range('E34').Value = '  -0.89%  '

$ws.Range('D36').Value = "'0.9995"

$ws.Range('E37').Value = '  -0.42%  '

$ws.Range('D38').Value = "'0.02028"
$ws.Range('E38').Value = '  -0.14%  '

$ws.Range('D39').Value = "'2.667"
$ws.Range('E39').Value = '  -0.76%  '

$ws.Range('D40').Value = "'110.90"
$ws.Range('E40').Value = '  +2.24%  '

$ws.Range('D41').Value = "'2.014"
$ws.Range('E41').Value = '  -1.63%  '

$ws.Range('D42').Value = "'0.4459"
$ws.Range('E42').Value = '  +6.08%  '

$ws.Range('D43').Value = "'0.8645"
$ws.Range('E43').Value = '  -1.03%  '

$ws.Range('E44').Value = '  -0.74%  '

$ws.Range('D45').Value = "'1.000"
$ws.Range('E45').Value = '  +0.05%  '

$ws.Range('D46').Value = "'68.10"
$ws.Range('E46').Value = '  +0.98%  '

$ws.Range('D47').Value = "'7.348"
$ws.Range('E47').Value = '  +2.56%  '

$ws.Range('D48').Value = "'48.96"
$ws.Range('E48').Value = '  -7.26%  '

$ws.Range('D49').Value = "'9.294"
$ws.Range('E49').Value = '  +1.16%  '

$ws.Range('E50').Value = '  +2.73%  '

$ws.Range('D51').Value = "'34.82"
$ws.Range('E51').Value = '  +0.43%  '
